# "Alterações no click criterio"
# Populate Sheet1 with the email/identifier header row and the first data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Endereço de e-mail"
$ws.Range("B1").Value = "`nIdentificador`n"

# Data row
$ws.Range("A2").Value = "lutisto@gmail.com"
$ws.Range("B2").Value = 50066528

# B1 keeps the existing wrap-text style; wrap + row height so the
# two-line header (blank / Identificador / blank) is fully visible.
$ws.Range("B1").WrapText = $true
$ws.Rows.Item(1).RowHeight = 43.2
